$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update existing rows 2-4: the "nutritionLabel" column (F) value
#    changes from "waffles" to "placeholder" for every existing item.
# ------------------------------------------------------------------
$ws.Range("F2").Value = "placeholder"
$ws.Range("F3").Value = "placeholder"
$ws.Range("F4").Value = "placeholder"

# ------------------------------------------------------------------
# 2. Add the new rows (5-11) with their values.
# ------------------------------------------------------------------

# Row 5 - Sausage Breakfast Sandwich
$ws.Range("A5").Value = "Sausage Breakfast Sandwich"
$ws.Range("B5").Value = "Muffin / Cooked Sausage Round / Cooked Egg Round / Pasteurized Process Cheddar Cheese"
$ws.Range("C5").Value = "Soy, wheat, milk, eggs."
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "Sausage_Breakfast_Sandwich"

# Row 6 - Bacon Breakfast Sandwich
$ws.Range("A6").Value = "Bacon Breakfast Sandwich"
$ws.Range("B6").Value = "English Muffin / Egg Patty / Processed Cheddar Cheese / Bacon."
$ws.Range("C6").Value = "Soy, wheat, milk, eggs."
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "Bacon_Breakfast_Sandwich"

# Row 7 - Egg and Cheese Breakfast Sandwich
$ws.Range("A7").Value = "Egg and Cheese Breakfast Sandwich"
$ws.Range("B7").Value = "English Muffin / Egg Patty / Processed Cheddar Cheese"
$ws.Range("C7").Value = "Soy, wheat, milk, eggs."
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "VEG"
$ws.Range("F7").Value = "Egg_and_Cheese_Breakfast_Sandwich"

# Row 8 - Tofu Scramble Roll
$ws.Range("A8").Value = "Tofu Scramble Roll"
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "VGN"
$ws.Range("F8").Value = "placeholder"

# Row 9 - Beef Sausage Roll
$ws.Range("A9").Value = "Beef Sausage Roll"
$ws.Range("D9").Value = "NA"
$ws.Range("F9").Value = "placeholder"

# Row 10 - Spicy Chicken Sausage Roll (no values besides A, F gets styled but empty)
$ws.Range("A10").Value = "Spicy Chicken Sausage Roll "

# Row 11 - Breakfast Waffles
$ws.Range("A11").Value = "Breakfast Waffles"
$ws.Range("B11").Value = "Belgian-style Waffles / Whipped Cream / Seasonal Fruit Compote or Maple Syrup "
$ws.Range("C11").Value = "Wheat, milk, eggs, "
$ws.Range("D11").Value = "NA"
$ws.Range("F11").Value = "placeholder"

# ------------------------------------------------------------------
# 3. Row heights for the newly inserted rows.
# ------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 30

# ------------------------------------------------------------------
# 4. Formatting.
# ------------------------------------------------------------------

# wrap text on the "Ingredients" cells for the new sandwich rows
$ws.Range("B5").WrapText = $true
$ws.Range("B6").WrapText = $true

# centered bold-ish icon column (ItemName) for the sandwich rows
$ws.Range("A5:A8").HorizontalAlignment = -4108
$ws.Range("A5:A8").VerticalAlignment = -4108
$ws.Range("A5:A8").Font.Color = 0

# nutritionLabel "placeholder" cells: green fill + thin white border
$nutritionCells = $ws.Range("F2,F3,F4,F8,F9,F10,F11")
$nutritionCells.Font.Color = 0
$nutritionCells.Interior.Color = 11854022
$nutritionCells.Interior.PatternColor = 11854022
$nutritionCells.Borders.Color = 16777215
$nutritionCells.Borders.LineStyle = 1
$nutritionCells.Borders.Weight = 2

# ------------------------------------------------------------------
# 5. Resize the Excel table (ListObject) to cover the new rows.
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G11"))

# ------------------------------------------------------------------
# 6. Selection, matching the authored workbook state.
# ------------------------------------------------------------------
$ws.Range("A10").Select()
